$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.00264
$ws.Range("E2").Value = -0.0745
$ws.Range("F2").Value = 0.108
$ws.Range("G2").Value = 0.1460353679406731
$ws.Range("H2").Value = 0.1460353679406731
$ws.Range("I2").Value = 0.09966936956483066
$ws.Range("J2").Value = 0.06366107437274614
$ws.Range("K2").Value = 96.09999999999999
$ws.Range("L2").Value = 0.05482030804335425
$ws.Range("M2").Value = 4.72
$ws.Range("N2").Value = 0.003134338269473404
$ws.Range("O2").Value = 0.04911550468262227
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 4.72
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 555
$ws.Range("V2").Value = 0.3685503685503685
$ws.Range("W2").Value = 0.04666181111920369
$ws.Range("X2").Value = 0.09353095774912788
$ws.Range("Y2").Value = -0.0468691466299242
$ws.Range("Z2").Value = 0.3728162726189691
$ws.Range("AA2").Value = 0.02373388445856619
$ws.Range("AB2").Value = 0.04529620934707808
$ws.Range("AC2").Value = -0.02156232488851189
$ws.Range("AD2").Value = 3079.7
$ws.Range("AE2").Value = 2.947975764259178
$ws.Range("AF2").Value = 3082.647975764259
$ws.Range("AG2").Value = 2527.647975764259
$ws.Range("AH2").Value = 0.6718133910871489
$ws.Range("AI2").Value = 0.5742693464396641
$ws.Range("AJ2").Value = 0.6266562319208144
$ws.Range("AK2").Value = 0.5251766668769962
$ws.Range("AL2").Value = 26.8
$ws.Range("AM2").Value = 26.8
$ws.Range("AN2").Value = 18.15753788102116
$ws.Range("AO2").Value = 6.440298507462686
$ws.Range("AP2").Value = 14.90270606546937
$ws.Range("AQ2").Value = 6.440298507462686

$ws.Range("D3").Value = -0.00264
$ws.Range("E3").Value = -0.0745
$ws.Range("F3").Value = 0.108
$ws.Range("G3").Value = 0.1460353679406731
$ws.Range("H3").Value = 0.1460353679406731
$ws.Range("I3").Value = 0.09966936956483066
$ws.Range("J3").Value = 0.06366107437274614
$ws.Range("K3").Value = 96.09999999999999
$ws.Range("L3").Value = 0.05482030804335425
$ws.Range("M3").Value = 4.72
$ws.Range("N3").Value = 0.003134338269473404
$ws.Range("O3").Value = 0.04911550468262227
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 4.72
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 555
$ws.Range("V3").Value = 0.3685503685503685
$ws.Range("W3").Value = 0.04666181111920369
$ws.Range("X3").Value = 0.09353095774912788
$ws.Range("Y3").Value = -0.0468691466299242
$ws.Range("Z3").Value = 0.3728162726189691
$ws.Range("AA3").Value = 0.02373388445856619
$ws.Range("AB3").Value = 0.04529620934707808
$ws.Range("AC3").Value = -0.02156232488851189
$ws.Range("AD3").Value = 3079.7
$ws.Range("AE3").Value = 2.947975764259178
$ws.Range("AF3").Value = 3082.647975764259
$ws.Range("AG3").Value = 2527.647975764259
$ws.Range("AH3").Value = 0.6718133910871489
$ws.Range("AI3").Value = 0.5742693464396641
$ws.Range("AJ3").Value = 0.6266562319208144
$ws.Range("AK3").Value = 0.5251766668769962
$ws.Range("AL3").Value = 26.8
$ws.Range("AM3").Value = 26.8
$ws.Range("AN3").Value = 18.15753788102116
$ws.Range("AO3").Value = 6.440298507462686
$ws.Range("AP3").Value = 14.90270606546937
$ws.Range("AQ3").Value = 6.440298507462686
